# Se procesan de nuevo los datos con las nuevas dimensiones curadas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (subespecie-ganaderia-descripcion) : dimension -> measure
$ws.Range("B2").Value = "iaest-measure:subespecie-ganaderia-descripcion"
$ws.Range("B3").Value = "medida"
$ws.Range("B4").Value = "xsd:int"
$ws.Range("B5").ClearContents()

# Column F (especie-ganaderia-descripcion) : dimension -> measure
$ws.Range("F2").Value = "iaest-measure:especie-ganaderia-descripcion"
$ws.Range("F3").Value = "medida"
$ws.Range("F4").Value = "xsd:int"
$ws.Range("F5").ClearContents()

# Column J (animales) : dimension -> measure
$ws.Range("J2").Value = "iaest-measure:animales"
$ws.Range("J3").Value = "medida"
$ws.Range("J4").Value = "xsd:int"
$ws.Range("J5").ClearContents()

# Column L (municipio-nombre) : measure -> dimension
$ws.Range("L2").Value = "sdmx-dimension:refArea"
$ws.Range("L3").Value = "dim"
$ws.Range("L4").Value = "URI-Municipio"
